$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginData")

# Update data values: keep header row, replace rows 2-3 with new data,
# and remove the now-unused rows 4-6.
$ws.Range("A2").Value = "mngr353180"
$ws.Range("B2").Value = "nerynYt"
$ws.Range("A3").Value = "mngr353180"
$ws.Range("B3").Value = "nerynYt"

$ws.Rows("4:6").Delete()

$ws.Range("A3").Select()
